$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.890575333333333
$ws.Range("H2").Value = 5.671726
$ws.Range("I2").Value = 0.006346320422088561
$ws.Range("J2").Value = 0.00634632042208856
$ws.Range("M2").Value = 11.319211
$ws.Range("N2").Value = 33.957633
$ws.Range("O2").Value = 0.09922284194232082
$ws.Range("P2").Value = 0.09922284194232082
$ws.Range("Q2").Value = 21.39982110939533
$ws.Range("R2").Value = 192.598389984558
$ws.Range("S2").Value = 0.000629699948156216
$ws.Range("T2").Value = 0.000629699948156216
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.890575333333333
$ws.Range("H3").Value = 5.671726
$ws.Range("I3").Value = 0.006346320422088561
$ws.Range("J3").Value = 0.00634632042208856
$ws.Range("O3").Value = 0.3843080175847637
$ws.Range("P3").Value = 0.3843080175847637
$ws.Range("Q3").Value = 82.885378671184
$ws.Range("R3").Value = 745.968408040656
$ws.Range("S3").Value = 0.002438941820370556
$ws.Range("T3").Value = 0.002438941820370555
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.890575333333333
$ws.Range("H4").Value = 5.671726
$ws.Range("I4").Value = 0.006346320422088561
$ws.Range("J4").Value = 0.00634632042208856
$ws.Range("O4").Value = 0.5164691404729155
$ws.Range("P4").Value = 0.5164691404729155
$ws.Range("Q4").Value = 111.3891418376062
$ws.Range("R4").Value = 1002.502276538456
$ws.Range("S4").Value = 0.003277678653561789
$ws.Range("T4").Value = 0.003277678653561789
$ws.Range("I5").Value = 0.8887896079640043
$ws.Range("J5").Value = 0.8887896079640044
$ws.Range("M5").Value = 11.319211
$ws.Range("N5").Value = 33.957633
$ws.Range("O5").Value = 0.09922284194232082
$ws.Range("P5").Value = 0.09922284194232082
$ws.Range("Q5").Value = 2997.002569885982
$ws.Range("R5").Value = 26973.02312897384
$ws.Range("S5").Value = 0.08818823079098968
$ws.Range("T5").Value = 0.0881882307909897
$ws.Range("I6").Value = 0.8887896079640043
$ws.Range("J6").Value = 0.8887896079640044
$ws.Range("O6").Value = 0.3843080175847637
$ws.Range("P6").Value = 0.3843080175847637
$ws.Range("S6").Value = 0.3415689722865858
$ws.Range("T6").Value = 0.3415689722865858
$ws.Range("I7").Value = 0.8887896079640043
$ws.Range("J7").Value = 0.8887896079640044
$ws.Range("O7").Value = 0.5164691404729155
$ws.Range("P7").Value = 0.5164691404729155
$ws.Range("S7").Value = 0.4590324048864288
$ws.Range("T7").Value = 0.4590324048864289
$ws.Range("I8").Value = 0.104864071613907
$ws.Range("J8").Value = 0.104864071613907
$ws.Range("M8").Value = 11.319211
$ws.Range("N8").Value = 33.957633
$ws.Range("O8").Value = 0.09922284194232082
$ws.Range("P8").Value = 0.09922284194232082
$ws.Range("Q8").Value = 353.602122819054
$ws.Range("R8").Value = 3182.419105371486
$ws.Range("S8").Value = 0.0104049112031749
$ws.Range("T8").Value = 0.0104049112031749
$ws.Range("I9").Value = 0.104864071613907
$ws.Range("J9").Value = 0.104864071613907
$ws.Range("O9").Value = 0.3843080175847637
$ws.Range("P9").Value = 0.3843080175847637
$ws.Range("S9").Value = 0.04030010347780727
$ws.Range("T9").Value = 0.04030010347780727
$ws.Range("I10").Value = 0.104864071613907
$ws.Range("J10").Value = 0.104864071613907
$ws.Range("O10").Value = 0.5164691404729155
$ws.Range("P10").Value = 0.5164691404729155
$ws.Range("S10").Value = 0.05415905693292478
$ws.Range("T10").Value = 0.05415905693292478
